$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $result = $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. Activation date
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2025"

# 2. Objectives (PT)
Replace-Text `
    "Capacitar os alunos para a aplicação de conhecimentos da engenharia química na solução de problemas que se apresentam na implantação e otimização de processos biotecnológicos, com ênfase em: 1) esterilização de equipamentos, meios e ar e; 2) cinética e estequiometria do crescimento microbiano e da formação de produtos." `
    "Desenvolver nos discentes as competências e habilidades necessárias para a aplicação de conhecimentos científicos, tecnológicos e de engenharia na concepção, projeto, instalação, otimização, supervisão e avaliação crítica da operação de bioprocessos, com ênfase em: 1) esterilização de equipamentos, meios e ar e; 2) cinética e estequiometria do crescimento microbiano e da formação de produtos."

# 3. Objectives (EN) - currently an empty italic run; find the paragraph right after the PT objectives
# and set its text (preserves the existing <w:i/> run formatting).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "" -and $p.Range.Font.Italic -eq -1) {
        $p.Range.Text = "To develop in students the competencies and skills necessary for the application of scientific, technological, and engineering knowledge in the conception, design, installation, optimization, supervision, and critical evaluation of the operation of bioprocesses, with emphasis on: 1) sterilization of equipment, media, and air; and 2) kinetics and stoichiometry of microbial growth and product formation."
        break
    }
}

# 4. Short program (PT)
Replace-Text `
    "Processo biotecnológico genérico, esterilização de equipamentos, esterilização de meios por aquecimento com vapor, esterilização de ar por filtração, cinética e estequiometria do crescimento microbiano e da formação de produtos." `
    "Processo biotecnológico genérico, esterilização de equipamentos, esterilização de meios por aquecimento com vapor, esterilização por filtração, cinética e estequiometria do crescimento microbiano e da formação de produtos."

# 5. Short program (EN)
Replace-Text `
    "Generic biotechnological process, equipment sterilization, media sterilization by heating, air sterilization by filtration, kinetics and stoichiometry of microbial growth and products formation." `
    "Generic biotechnological process, equipment sterilization, media sterilization by steam heating, sterilization by filtration, kinetics and stoichiometry of microbial growth and product formation."

# 6. Full program (PT)
Replace-Text `
    "1. Processo biotecnológico genérico: representação esquemática; descrição das etapas principais.2. Esterilização de equipamentos: terminologia; esterilização por agentes físicos; esterilização por agentes químicos.3. Esterilização de meios por aquecimento com vapor: cinética da destruição térmica de microrganismos; destruição de nutrientes do meio; cálculo do tempo de esterilização por processo descontínuo; cálculo do tempo de esterilização por processo contínuo.4. Esterilização de ar por filtração: aerossóis microbianos; amostradores; dimensionamento de filtros fibrosos; dimensionamento de filtros de membranas.5. Cinética e estequiometria do crescimento microbiano e da formação de produtos: velocidades de transformação e fatores de conversão; classificação dos processos fermentativos em função das velocidades de crescimento celular e formação de produtos; influência da concentração do substrato sobre a velocidade de crescimento celular; estequiometria do crescimento microbiano e da formação de produtos." `
    "1. Processo biotecnológico genérico: representação esquemática; descrição das etapas principais. 2. Esterilização de equipamentos: terminologia; esterilização por agentes físicos; esterilização por agentes químicos. 3. Esterilização de meios por aquecimento com vapor: cinética da destruição térmica de microrganismos; destruição de nutrientes do meio; cálculo do tempo de esterilização por processo descontínuo; dimensionamento de sistemas de esterilização por processo contínuo. 4. Esterilização por filtração: aerossóis microbianos; amostradores; dimensionamento de filtros; esterilização de meios. 5. Cinética e estequiometria do crescimento microbiano e da formação de produtos: velocidades de transformação e fatores de conversão; classificação dos processos fermentativos em função das velocidades de crescimento celular e formação de produtos; influência da concentração do substrato sobre a velocidade de crescimento celular; estequiometria do crescimento microbiano e da formação de produtos."

# 7. Full program (EN)
Replace-Text `
    "1.Generic biotechnological process: schematic representation; description of the main phases.2.Equipment sterilization: terminology, sterilization by physical agents, sterilization by chemical agents.3.Media sterilization by heating: kinetics of thermic destruction of microorganisms; destruction of medium nutrients; calculation of sterilization time by discontinuous process; calculation of sterilization time by continuous process.4.Air sterilization by filtration: microbial aerosols; air samplers; dimensioning of fibrous filters; dimensioning of membrane filters.5.Kinetics and stoichiometry of microbial growth and products formation: definition of velocities of transformation and conversion factors; classification of fermentations as a function of the velocities of cell growth and products formation; influence of substrate concentration on cell growth; stoichiometry of microbial growth and products formation." `
    "1.Generic biotechnological process: schematic representation; description of the main stages.2.Equipment sterilization: terminology; sterilization by physical agents; sterilization by chemical agents.3.Media sterilization by steam heating: kinetics of thermal destruction of microorganisms; destruction of nutrient media; calculation of sterilization time for batch processes; design of sterilization systems for continuous processes.4.Sterilization by filtration: microbial aerosols; air samplers; filter sizing; media sterilization.5.Kinetics and stoichiometry of microbial growth and product formation: transformation rates and conversion factors; classification of fermentative processes based on cell growth and product formation rates; influence of substrate concentration on cell growth rate; stoichiometry of microbial growth and product formation."

# 8. Método
Replace-Text `
    "Os alunos serão avaliados formalmente por duas provas escritas (P1 e P2), sendo a segunda prova (P2) com peso 2." `
    "A avaliação do aprendizado será feita pela aplicação de duas provas escritas, e através de trabalhos desenvolvidos pelos discentes (estes poderão incluir estudos dirigidos, análises de artigos, resolução de problemas práticos, entre outros)."

# 9. Critério
Replace-Text `
    "A nota final (NF) será calculada como: NF=(P1+(P2×2))/3. Serão aprovados os alunos que obtiverem NF maior ou igual 5,0." `
    "A nota final (NF) será composta pelas médias M1  e M2,calculadas conforme segue:M1=P1+a1×T1M2=P2+a2×T2Em que:-P1 e P2 são as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).-T1 e T2 são as notas médias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.-a1 e a2 são os fatores multiplicadores das notas médias dos trabalhos, a serem definidos pelo docente antes do início de cada turma com base nas atividades específicas a serem propostas. Os valores serão ≥0,1, sendo informados aos alunos no início do semestre. Em todos os casos, os valores máximos para M1 e M2 serão “dez”, sendo desconsideradas pontuações superiores.O cálculo de NF será feito conforme segue:NF=(M1+2×M2)/3Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."

# 10. Norma de recuperação - merges two <w:t> runs separated by a <w:br/> into one <w:t>
$vbreak = [char]11
$oldRecup = "Será oferecido um programa de recuperação avaliado por uma prova escrita final (PR)." + $vbreak + "A média de recuperação (MR) será calculada como: MR=(NF+PR)/2. Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$newRecup = "Será oferecido um programa de recuperação, sendo este avaliado por uma prova escrita final (PR). A média de recuperação (MR) será calculada conforme segue: MR=(NF+PR)/2Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
Replace-Text $oldRecup $newRecup

# 11. Bibliography
Replace-Text `
    "BAILEY, J.E., OLLIS D.F. Biochemical Engineering Fundamentals. 2nd edition, New York: McGraw Hill, 1986. ISBN: 978-0070032125.BORZANI, W., SCHMIDELL, W., LIMA, U.A., AQUARONE, E. Biotecnologia Industrial Fundamentos (Vol 1). São Paulo: Edgard Blucher Ltda, 2001.m ISBN: 9788521202783.DORAN P.M. Bioprocess Engineering Principles, 1st edition, San Diego: Academic Press, 1995. ISBN: 978-0080528120.KATOH, S., HORIUCHI, J., YOSHIDA, F. Biochemical Engineering: A Textbook for Engineers, Chemists and Biologists, 2nd, Completely Revised and Enlarged Edition. Weinheim/Germany: Wiley-VCH, 2015. ISBN: 978-3527338047.SCHMIDELL, W., LIMA, U.A., AQUARONE, E., BORZANI, W. Biotecnologia Industrial Engenharia Bioquímica (Vol 2), São Paulo: Edgard Blucher Ltda, 2001. ISBN: 9788521202790." `
    "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. Biotecnologia Industrial. Volume 1: Fundamentos. 2ª Edição. São Paulo: Blucher, 2020. ISBN 978-85-212-1897-5 (e-Book); 978-85-212-1898-2 (Impresso).ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. 2ª Edição. São Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).BORZANI, W. Processo Biotecnológico Industrial Genérico. In: BORZANI, W.; SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E. Biotecnologia Industrial. Volume 1: Fundamentos. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0278-3.DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706"

Write-Output "Done"
